$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "62.559.28"
$ws.Cells.Item(2, 5).Value = "  +3.13%  "
$ws.Cells.Item(3, 4).Value = "2.439.50"
$ws.Cells.Item(3, 5).Value = "  +1.60%  "
$ws.Cells.Item(4, 5).Value = "  -0.14%  "
$ws.Cells.Item(5, 4).Formula = "'577.22"
$ws.Cells.Item(5, 5).Value = "  +2.22%  "
$ws.Cells.Item(6, 4).Formula = "'145.13"
$ws.Cells.Item(6, 5).Value = "  +2.86%  "
$ws.Cells.Item(7, 5).Value = "  +0.07%  "
$ws.Cells.Item(8, 5).Value = "  +0.11%  "
$ws.Cells.Item(9, 4).Value = "2.437.81"
$ws.Cells.Item(9, 5).Value = "  +1.23%  "
$ws.Cells.Item(10, 5).Value = "  +1.87%  "
$ws.Cells.Item(11, 5).Value = "  +1.03%  "
$ws.Cells.Item(12, 5).Value = "  +0.53%  "
$ws.Cells.Item(14, 4).Formula = "'28.36"
$ws.Cells.Item(14, 5).Value = "  +8.74%  "
$ws.Cells.Item(15, 5).Value = "  +5.11%  "
$ws.Cells.Item(16, 4).Value = "2.881.38"
$ws.Cells.Item(16, 5).Value = "  +2.50%  "
$ws.Cells.Item(17, 4).Value = "62.457.33"
$ws.Cells.Item(17, 5).Value = "  +3.22%  "
$ws.Cells.Item(18, 4).Value = "0.0₆0914"
$ws.Cells.Item(18, 5).Value = "  +218.75%  "
$ws.Cells.Item(19, 4).Value = "2.436.93"
$ws.Cells.Item(19, 5).Value = "  +1.46%  "
$ws.Cells.Item(20, 5).Value = "  -3.61%  "
$ws.Cells.Item(21, 4).Formula = "'10.87"
$ws.Cells.Item(21, 5).Value = "  +2.26%  "
$ws.Cells.Item(22, 4).Formula = "'325.47"
$ws.Cells.Item(22, 5).Value = "  +0.39%  "
$ws.Cells.Item(23, 5).Value = "  +1.06%  "
$ws.Cells.Item(24, 4).Formula = "'2.02"
$ws.Cells.Item(24, 5).Value = "  +9.80%  "
$ws.Cells.Item(25, 5).Value = "  +0.02%  "
$ws.Cells.Item(26, 4).Formula = "'65.25"
$ws.Cells.Item(26, 5).Value = "  +0.25%  "
$ws.Cells.Item(27, 4).Formula = "'633.99"
$ws.Cells.Item(27, 5).Value = "  +11.46%  "
$ws.Cells.Item(28, 5).Value = "  +13.74%  "
$ws.Cells.Item(29, 4).Formula = "'8.49"
$ws.Cells.Item(29, 5).Value = "  +5.50%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0975"
$ws.Cells.Item(30, 5).Value = "  +4.38%  "
$ws.Cells.Item(32, 4).Formula = "'8.18"
$ws.Cells.Item(32, 5).Value = "  +1.11%  "
$ws.Cells.Item(33, 5).Value = "  +5.11%  "
$ws.Cells.Item(34, 5).Value = "  +5.77%  "
$ws.Cells.Item(35, 5).Value = "  +2.29%  "
$ws.Cells.Item(36, 5).Value = "  +1.64%  "
$ws.Cells.Item(37, 5).Value = "  +0.10%  "
$ws.Cells.Item(38, 5).Value = "  +3.13%  "
$ws.Cells.Item(39, 4).Formula = "'153.00"
$ws.Cells.Item(39, 5).Value = "  -0.17%  "
$ws.Cells.Item(40, 5).Value = "  +5.91%  "
$ws.Cells.Item(41, 5).Value = "  +0.01%  "
$ws.Cells.Item(42, 4).Formula = "'18.54"
$ws.Cells.Item(42, 5).Value = "  +1.24%  "
$ws.Cells.Item(43, 5).Value = "  +8.95%  "
$ws.Cells.Item(44, 5).Value = "  +4.43%  "
$ws.Cells.Item(45, 4).Formula = "'42.42"
$ws.Cells.Item(45, 5).Value = "  +1.62%  "
$ws.Cells.Item(46, 4).Formula = "'0.999"
$ws.Cells.Item(46, 5).Value = "  +0.06%  "
$ws.Cells.Item(47, 5).Value = "  +28.03%  "
$ws.Cells.Item(48, 4).Formula = "'144.05"
$ws.Cells.Item(48, 5).Value = "  +1.27%  "
$ws.Cells.Item(49, 5).Value = "  +0.81%  "
$ws.Cells.Item(50, 4).Formula = "'20.43"
$ws.Cells.Item(50, 5).Value = "  +5.66%  "
$ws.Cells.Item(51, 5).Value = "  +1.57%  "
